$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "After Adaptation" data block (rows 17-25, columns B:AA).
# AB (row total) and AC (percent change) are formulas that recompute automatically.
$row17 = New-Object 'object[,]' 1,26
$row17[0,0] = 17783.939999999999
$row17[0,1] = 578.37
$row17[0,2] = 2230.27
$row17[0,3] = 9794.65
$row17[0,4] = 46771.66
$row17[0,5] = 8794.11
$row17[0,6] = 7365.9699999999903
$row17[0,7] = 12380.91
$row17[0,8] = 1650.61
$row17[0,9] = 6989.36
$row17[0,10] = 2567.33
$row17[0,11] = 12044.66
$row17[0,12] = 6740.08
$row17[0,13] = 1228.75
$row17[0,14] = 1160.3599999999999
$row17[0,15] = 15965.65
$row17[0,16] = 3071.76
$row17[0,17] = 9215.24
$row17[0,18] = 4098.83
$row17[0,19] = 7592.7
$row17[0,20] = 11928.53
$row17[0,21] = 1414.57
$row17[0,22] = 21664.22
$row17[0,23] = 8980.66
$row17[0,24] = 2751.44
$row17[0,25] = 43617.27
$ws.Range("B17:AA17").Value = $row17

$row18 = New-Object 'object[,]' 1,26
$row18[0,0] = 16320.25
$row18[0,1] = 542.07999999999902
$row18[0,2] = 2096.89
$row18[0,3] = 8876.33
$row18[0,4] = 40147.31
$row18[0,5] = 7988.9
$row18[0,6] = 6758.3899999999903
$row18[0,7] = 11044.7599999999
$row18[0,8] = 1529.76
$row18[0,9] = 6283.71
$row18[0,10] = 2298.3599999999901
$row18[0,11] = 11160.619999999901
$row18[0,12] = 6094.82
$row18[0,13] = 1134.08
$row18[0,14] = 1110.55
$row18[0,15] = 15035.14
$row18[0,16] = 2771.88
$row18[0,17] = 8438.7199999999993
$row18[0,18] = 3898.69
$row18[0,19] = 6991.86
$row18[0,20] = 11028.11
$row18[0,21] = 1307.73
$row18[0,22] = 15456.64
$row18[0,23] = 8007.1399999999903
$row18[0,24] = 2523.04
$row18[0,25] = 41009.629999999997
$ws.Range("B18:AA18").Value = $row18

$row19 = New-Object 'object[,]' 1,26
$row19[0,0] = 18956.77
$row19[0,1] = 593.17999999999995
$row19[0,2] = 2294.0299999999902
$row19[0,3] = 10531.06
$row19[0,4] = 54368.539999999899
$row19[0,5] = 9145.67
$row19[0,6] = 7708.99
$row19[0,7] = 13172.13
$row19[0,8] = 1709.63
$row19[0,9] = 7348.37
$row19[0,10] = 2675.73
$row19[0,11] = 12840.65
$row19[0,12] = 7043.84
$row19[0,13] = 1271.96
$row19[0,14] = 1195.76
$row19[0,15] = 16580.52
$row19[0,16] = 3202.19
$row19[0,17] = 9549.84
$row19[0,18] = 4208.62
$row19[0,19] = 7864.51
$row19[0,20] = 12500.29
$row19[0,21] = 1480.46
$row19[0,22] = 26246
$row19[0,23] = 9900.3799999999992
$row19[0,24] = 2859.47
$row19[0,25] = 45785.84
$ws.Range("B19:AA19").Value = $row19

$row20 = New-Object 'object[,]' 1,26
$row20[0,0] = 18196.21
$row20[0,1] = 589.42999999999995
$row20[0,2] = 2277.84
$row20[0,3] = 10125.33
$row20[0,4] = 47868.03
$row20[0,5] = 8912.36
$row20[0,6] = 7519.6
$row20[0,7] = 12753.35
$row20[0,8] = 1683.3899999999901
$row20[0,9] = 7126.8799999999901
$row20[0,10] = 2600.3199999999902
$row20[0,11] = 12420.13
$row20[0,12] = 6807.77
$row20[0,13] = 1247.17
$row20[0,14] = 1188.32
$row20[0,15] = 16389.419999999998
$row20[0,16] = 3128.16
$row20[0,17] = 9347.5400000000009
$row20[0,18] = 4192.7700000000004
$row20[0,19] = 7706.7199999999903
$row20[0,20] = 12173.48
$row20[0,21] = 1449.57
$row20[0,22] = 23205.14
$row20[0,23] = 9384.4199999999892
$row20[0,24] = 2796.6
$row20[0,25] = 44468.51
$ws.Range("B20:AA20").Value = $row20

$row21 = New-Object 'object[,]' 1,26
$row21[0,0] = 16258.2
$row21[0,1] = 553.97
$row21[0,2] = 2133.5
$row21[0,3] = 9122.66
$row21[0,4] = 40822.699999999997
$row21[0,5] = 7951.3099999999904
$row21[0,6] = 6749.9699999999903
$row21[0,7] = 11073.81
$row21[0,8] = 1530.81
$row21[0,9] = 6209.98
$row21[0,10] = 2270.71
$row21[0,11] = 11433.4
$row21[0,12] = 6086.98
$row21[0,13] = 1159.3699999999999
$row21[0,14] = 1137.4000000000001
$row21[0,15] = 15398.91
$row21[0,16] = 2768.5
$row21[0,17] = 8415.91
$row21[0,18] = 3965.2
$row21[0,19] = 6961.43
$row21[0,20] = 11087.6899999999
$row21[0,21] = 1316.25
$row21[0,22] = 19849.949999999899
$row21[0,23] = 8391.64
$row21[0,24] = 2496.5299999999902
$row21[0,25] = 41094.699999999997
$ws.Range("B21:AA21").Value = $row21

$row22 = New-Object 'object[,]' 1,26
$row22[0,0] = 19573.72
$row22[0,1] = 607.16999999999996
$row22[0,2] = 2342.91
$row22[0,3] = 10897.11
$row22[0,4] = 56573.2599999999
$row22[0,5] = 9360.02
$row22[0,6] = 7892.67
$row22[0,7] = 13646.74
$row22[0,8] = 1755.77
$row22[0,9] = 7552.24
$row22[0,10] = 2733.76
$row22[0,11] = 13287.03
$row22[0,12] = 7185.13
$row22[0,13] = 1302.52
$row22[0,14] = 1234.8899999999901
$row22[0,15] = 17005.2
$row22[0,16] = 3286.26
$row22[0,17] = 9798.82
$row22[0,18] = 4325.95
$row22[0,19] = 8044.01
$row22[0,20] = 12793.51
$row22[0,21] = 1520.85
$row22[0,22] = 27766.85
$row22[0,23] = 10269.82
$row22[0,24] = 2928.18
$row22[0,25] = 47285.59
$ws.Range("B22:AA22").Value = $row22

$row23 = New-Object 'object[,]' 1,26
$row23[0,0] = 18213.54
$row23[0,1] = 595.49
$row23[0,2] = 2295.73
$row23[0,3] = 9958.98
$row23[0,4] = 48333.63
$row23[0,5] = 8920.26
$row23[0,6] = 7486.53
$row23[0,7] = 12526.77
$row23[0,8] = 1689.9
$row23[0,9] = 7088.88
$row23[0,10] = 2603.9499999999998
$row23[0,11] = 12312.9199999999
$row23[0,12] = 6813.17
$row23[0,13] = 1249.98
$row23[0,14] = 1194.98999999999
$row23[0,15] = 16379.99
$row23[0,16] = 3131.31
$row23[0,17] = 9390.91
$row23[0,18] = 4236.1499999999996
$row23[0,19] = 7781.28
$row23[0,20] = 12152.65
$row23[0,21] = 1436.55
$row23[0,22] = 21881.67
$row23[0,23] = 9078.3799999999992
$row23[0,24] = 2801.22
$row23[0,25] = 44780.25
$ws.Range("B23:AA23").Value = $row23

$row24 = New-Object 'object[,]' 1,26
$row24[0,0] = 16267.869999999901
$row24[0,1] = 552.12
$row24[0,2] = 2126.7799999999902
$row24[0,3] = 8736.7199999999993
$row24[0,4] = 39772.839999999997
$row24[0,5] = 7788.46
$row24[0,6] = 6619.67
$row24[0,7] = 10617.52
$row24[0,8] = 1505.21
$row24[0,9] = 6027.92
$row24[0,10] = 2222.59
$row24[0,11] = 11099.45
$row24[0,12] = 5955.12
$row24[0,13] = 1137.1500000000001
$row24[0,14] = 1129.49
$row24[0,15] = 15258.95
$row24[0,16] = 2718.42
$row24[0,17] = 8273.3699999999899
$row24[0,18] = 3965.3999999999901
$row24[0,19] = 6946.15
$row24[0,20] = 10842.93
$row24[0,21] = 1290.78
$row24[0,22] = 14451.84
$row24[0,23] = 7783.68
$row24[0,24] = 2454.37
$row24[0,25] = 40895.82
$ws.Range("B24:AA24").Value = $row24

$row25 = New-Object 'object[,]' 1,26
$row25[0,0] = 19757.39
$row25[0,1] = 617.62
$row25[0,2] = 2385.27
$row25[0,3] = 10932.34
$row25[0,4] = 57311.3
$row25[0,5] = 9453.14
$row25[0,6] = 7945.9699999999903
$row25[0,7] = 13626.78
$row25[0,8] = 1781.49
$row25[0,9] = 7656.97
$row25[0,10] = 2762.73
$row25[0,11] = 13359.71
$row25[0,12] = 7262.51
$row25[0,13] = 1324.23
$row25[0,14] = 1260.19
$row25[0,15] = 17205.79
$row25[0,16] = 3318.41
$row25[0,17] = 9906.0299999999897
$row25[0,18] = 4414.57
$row25[0,19] = 8181.01
$row25[0,20] = 12927
$row25[0,21] = 1531.21999999999
$row25[0,22] = 26979.18
$row25[0,23] = 10222.11
$row25[0,24] = 2959.2799999999902
$row25[0,25] = 48132.88
$ws.Range("B25:AA25").Value = $row25


# Rows 18, 21 and 24 had their AB/AC formula cells re-entered directly,
# which detaches them from the shared-formula group used by the other rows.
$ws.Range("AB18").Formula = "=SUM(B18:AA18)"
$ws.Range("AC18").Formula = '=(AB18-$AB$2)/$AB$2'
$ws.Range("AB21").Formula = "=SUM(B21:AA21)"
$ws.Range("AC21").Formula = '=(AB21-$AB$2)/$AB$2'
$ws.Range("AB24").Formula = "=SUM(B24:AA24)"
$ws.Range("AC24").Formula = '=(AB24-$AB$2)/$AB$2'

# Restore the active selection that was saved with the workbook.
$ws.Range("A29").Select()
